$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.542.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.72%  "
$ws.Range("D3").Value = "'3.494.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.58%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'591.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.24%  "
$ws.Range("D6").Value = "'168.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.40%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'3.491.49"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.49%  "
$ws.Range("D9").Value = "'0.593"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.65%  "
$ws.Range("D10").Value = "'7.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "  +5.95%  "
$ws.Range("E12").Value = "  +3.01%  "
$ws.Range("D13").Value = "'4.095.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.42%  "
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").Value = "'28.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.22%  "
$ws.Range("D16").Value = "'66.570.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.74%  "
$ws.Range("E17").Value = "  +2.68%  "
$ws.Range("D18").Value = "'3.516.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.68%  "
$ws.Range("D19").Value = "'6.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("D20").Value = "'13.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.83%  "
$ws.Range("D21").Value = "'389.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.96%  "
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("D23").Value = "'72.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.56%  "
$ws.Range("E25").Value = "  +3.48%  "
$ws.Range("E26").Value = "  +5.09%  "
$ws.Range("D27").Value = "'10.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.21%  "
$ws.Range("D28").Value = "'0.179"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("D29").Value = "'0.996"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("E30").Value = "  +4.13%  "
$ws.Range("E31").Value = "  +4.14%  "
$ws.Range("E32").Value = "  +2.58%  "
$ws.Range("E33").Value = "  +3.55%  "
$ws.Range("E34").Value = "  +4.73%  "
$ws.Range("D35").Value = "'1.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.42%  "
$ws.Range("D36").Value = "'162.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.72%  "
$ws.Range("D37").Value = "'0.890"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.53%  "
$ws.Range("D38").Value = "'1.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.72%  "
$ws.Range("E39").Value = "  +4.89%  "
$ws.Range("D40").Value = "'0.0741"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.86%  "
$ws.Range("D41").Value = "'4.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.30%  "
$ws.Range("D42").Value = "'26.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.21%  "
$ws.Range("D43").Value = "'2.784.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("D44").Value = "'26.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.74%  "
$ws.Range("D45").Value = "'42.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("D46").Value = "'2.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.08%  "
$ws.Range("E47").Value = "  +1.93%  "
$ws.Range("D48").Value = "'344.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.77%  "
$ws.Range("E49").Value = "  +4.14%  "
$ws.Range("D50").Value = "'33.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.46%  "
$ws.Range("D51").Value = "'0.861"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.75%  "
